# Update cover sheet content and rename the data sheet in the remaining
# scenario file, matching the other already-updated scenario workbooks.

$wb = $excel.ActiveWorkbook

# --- Cover sheet: replace the generic template placeholders with the
#     actual values for this scenario file -------------------------------
$cover = $wb.Worksheets.Item("Cover")

$cover.Range("B19").Value = "Scenario"
$cover.Range("B23").Value = "Olexandr Balyk (UCC, olexandr.balyk@ucc.ie)"
$cover.Range("B20").Value = "Industry sector (IND)"
$cover.Range("B21").Value = "Specify combustion-related emission coefficients"
$cover.Range("B26").Value = "Olexandr Balyk (UCC, olexandr.balyk@ucc.ie)"

$cover.Activate() | Out-Null
$cover.Range("B23:D23").Select() | Out-Null

# --- Rename the data sheet from the generic "FLO_EMIS" to "Emissions" ---
$wb.Worksheets.Item("FLO_EMIS").Name = "Emissions"

# --- Drop the stale external reference to the old Cover Template file ---
$sources = $wb.LinkSources(1)
if ($sources) {
    foreach ($source in $sources) {
        $wb.BreakLink($source, 1)
    }
}
